# New crime data collected — update the weekly CompStat report:
#   - Volume/Number and report-week dates in the header
#   - Week-to-date / 28-day / YTD / 2-year crime-complaint figures (rows 14-33)
#   - Historical-perspective TOTAL-by-year row labels already match (no edit needed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/number and report week dates -----------------------
$ws.Range("A8").Value = "Volume 32   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/10/2025  Through  11/16/2025"

# --- Row 14 (Murder): C/D/F were text "0" and E was text "***.*"; they
#     now hold real numbers, so pull the numeric style from sibling cells
#     on the same row before writing the new values. ----------------------
$ws.Range("C14").NumberFormat = $ws.Range("G14").NumberFormat
$ws.Range("D14").NumberFormat = $ws.Range("G14").NumberFormat
$ws.Range("F14").NumberFormat = $ws.Range("G14").NumberFormat
$ws.Range("E14").NumberFormat = $ws.Range("H14").NumberFormat

$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = -50
$ws.Range("I14").Value = 26
$ws.Range("J14").Value = 46
$ws.Range("K14").Value = -43.478260869565
$ws.Range("L14").Value = -46.938775510204
$ws.Range("M14").Value = -43.478260869565
$ws.Range("N14").Value = -90.780141843971

# --- Row 15 (Rape) -------------------------------------------------------
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = -66.666666666666
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 170
$ws.Range("J15").Value = 134
$ws.Range("K15").Value = 26.865671641791
$ws.Range("L15").Value = 40.495867768595
$ws.Range("M15").Value = -10.526315789473
$ws.Range("N15").Value = -62.389380530973

# --- Row 16 (Robbery) -----------------------------------------------------
$ws.Range("C16").Value = 43
$ws.Range("D16").Value = 38
$ws.Range("E16").Value = 13.157894736842
$ws.Range("F16").Value = 136
$ws.Range("G16").Value = 171
$ws.Range("H16").Value = -20.467836257309
$ws.Range("I16").Value = 1531
$ws.Range("J16").Value = 1815
$ws.Range("K16").Value = -15.647382920110
$ws.Range("L16").Value = -12.862834376778
$ws.Range("M16").Value = -27.302943969610
$ws.Range("N16").Value = -81.962770970782

# --- Row 17 (Fel. Assault) -------------------------------------------------
$ws.Range("C17").Value = 47
$ws.Range("D17").Value = 66
$ws.Range("E17").Value = -28.787878787878
$ws.Range("F17").Value = 216
$ws.Range("G17").Value = 279
$ws.Range("H17").Value = -22.580645161290
$ws.Range("I17").Value = 2593
$ws.Range("J17").Value = 2919
$ws.Range("K17").Value = -11.168208290510
$ws.Range("L17").Value = -3.318419090231
$ws.Range("M17").Value = 52.619187757504
$ws.Range("N17").Value = -50.192086054552

# --- Row 18 (Burglary) ----------------------------------------------------
$ws.Range("C18").Value = 16
$ws.Range("D18").Value = 26
$ws.Range("E18").Value = -38.461538461538
$ws.Range("F18").Value = 97
$ws.Range("G18").Value = 130
$ws.Range("H18").Value = -25.384615384615
$ws.Range("I18").Value = 1306
$ws.Range("J18").Value = 1232
$ws.Range("K18").Value = 6.006493506493
$ws.Range("L18").Value = -2.391629297458
$ws.Range("M18").Value = 4.229848363926
$ws.Range("N18").Value = -86.919070512820

# --- Row 19 (Gr. Larceny) -------------------------------------------------
$ws.Range("D19").Value = 118
$ws.Range("E19").Value = -5.084745762711
$ws.Range("F19").Value = 458
$ws.Range("G19").Value = 505
$ws.Range("H19").Value = -9.306930693069
$ws.Range("I19").Value = 5549
$ws.Range("J19").Value = 5724
$ws.Range("K19").Value = -3.057302585604
$ws.Range("L19").Value = -3.813485872768
$ws.Range("M19").Value = 33.646435452793
$ws.Range("N19").Value = -42.917395329698

# --- Row 20 (G.L.A.) -------------------------------------------------------
$ws.Range("C20").Value = 18
$ws.Range("D20").Value = 19
$ws.Range("E20").Value = -5.263157894736
$ws.Range("F20").Value = 54
$ws.Range("G20").Value = 74
$ws.Range("H20").Value = -27.027027027027
$ws.Range("I20").Value = 769
$ws.Range("J20").Value = 849
$ws.Range("K20").Value = -9.422850412249
$ws.Range("L20").Value = -35.378151260504
$ws.Range("M20").Value = 49.031007751938
$ws.Range("N20").Value = -90.793726804740

# --- Row 21 (TOTAL) ---------------------------------------------------------
$ws.Range("C21").Value = 238
$ws.Range("D21").Value = 271
$ws.Range("E21").Value = -12.177121771217
$ws.Range("F21").Value = 974
$ws.Range("G21").Value = 1173
$ws.Range("H21").Value = -16.965046888320
$ws.Range("I21").Value = 11944
$ws.Range("J21").Value = 12719
$ws.Range("K21").Value = -6.093246324396
$ws.Range("L21").Value = -7.453897412056
$ws.Range("M21").Value = 19.895603292511
$ws.Range("N21").Value = -71.887209904439

# --- Row 22 (Transit) -------------------------------------------------------
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 33.333333333333
$ws.Range("G22").Value = 25
$ws.Range("H22").Value = -12
$ws.Range("I22").Value = 191
$ws.Range("J22").Value = 231
$ws.Range("K22").Value = -17.316017316017
$ws.Range("L22").Value = -30.03663003663
$ws.Range("M22").Value = -7.281553398058

# --- Row 23 (Housing) -------------------------------------------------------
$ws.Range("C23").Value = 31
$ws.Range("D23").Value = 27
$ws.Range("E23").Value = 14.814814814814
$ws.Range("F23").Value = 102
$ws.Range("G23").Value = 104
$ws.Range("H23").Value = -1.923076923076
$ws.Range("I23").Value = 1103
$ws.Range("J23").Value = 1167
$ws.Range("K23").Value = -5.484147386461
$ws.Range("L23").Value = -5.484147386461
$ws.Range("M23").Value = 46.092715231788

# --- Row 24 (Petit Larceny) ----------------------------------------------
$ws.Range("C24").Value = 267
$ws.Range("D24").Value = 304
$ws.Range("E24").Value = -12.171052631578
$ws.Range("F24").Value = 1010
$ws.Range("G24").Value = 1149
$ws.Range("H24").Value = -12.097476066144
$ws.Range("I24").Value = 12005
$ws.Range("J24").Value = 11855
$ws.Range("K24").Value = 1.265288907633
$ws.Range("L24").Value = -2.334851936218
$ws.Range("M24").Value = 43.343283582089

# --- Row 25 (Retail Theft) -------------------------------------------------
$ws.Range("C25").Value = 133
$ws.Range("D25").Value = 185
$ws.Range("E25").Value = -28.108108108108
$ws.Range("F25").Value = 492
$ws.Range("G25").Value = 679
$ws.Range("H25").Value = -27.540500736377
$ws.Range("I25").Value = 6232
$ws.Range("J25").Value = 6445
$ws.Range("K25").Value = -3.304887509697
$ws.Range("L25").Value = -8.177397966701

# --- Row 26 (Misd. Assault) -------------------------------------------------
$ws.Range("C26").Value = 88
$ws.Range("D26").Value = 80
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 333
$ws.Range("G26").Value = 399
$ws.Range("H26").Value = -16.541353383458
$ws.Range("I26").Value = 4250
$ws.Range("J26").Value = 4566
$ws.Range("K26").Value = -6.920718353044
$ws.Range("L26").Value = 4.064642507345
$ws.Range("M26").Value = -10.714285714285

# --- Row 27 (UCR Rape*) -----------------------------------------------------
$ws.Range("C27").Value = 3
$ws.Range("E27").Value = -25
$ws.Range("F27").Value = 18
$ws.Range("G27").Value = 16
$ws.Range("H27").Value = 12.5
$ws.Range("I27").Value = 207
$ws.Range("J27").Value = 207
$ws.Range("L27").Value = -7.589285714285

# --- Row 28 (Other Sex Crimes) ---------------------------------------------
$ws.Range("C28").Value = 11
$ws.Range("D28").Value = 11
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 46
$ws.Range("G28").Value = 45
$ws.Range("H28").Value = 2.222222222222
$ws.Range("I28").Value = 550
$ws.Range("J28").Value = 533
$ws.Range("K28").Value = 3.189493433395
$ws.Range("L28").Value = 6.589147286821

# --- Row 29 (Shooting Vic.) -------------------------------------------------
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = -80
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 18
$ws.Range("H29").Value = -72.222222222222
$ws.Range("I29").Value = 70
$ws.Range("J29").Value = 123
$ws.Range("K29").Value = -43.089430894308
$ws.Range("L29").Value = -46.564885496183
$ws.Range("M29").Value = -58.579881656804
$ws.Range("N29").Value = -90.112994350282

# --- Row 30 (Shooting Inc.) -------------------------------------------------
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 5
$ws.Range("E30").Value = -80
$ws.Range("F30").Value = 4
$ws.Range("G30").Value = 15
$ws.Range("H30").Value = -73.333333333333
$ws.Range("I30").Value = 59
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = -41
$ws.Range("L30").Value = -49.572649572649
$ws.Range("M30").Value = -59.589041095890
$ws.Range("N30").Value = -90.852713178294

# --- Row 31 (Hate Crimes) ---------------------------------------------------
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = -20
$ws.Range("I31").Value = 60
$ws.Range("J31").Value = 106
$ws.Range("K31").Value = -43.396226415094
$ws.Range("L31").Value = -40

# --- Row 33 (Traffic Fatalities) --------------------------------------------
$ws.Range("D33").Value = 1
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = -75
$ws.Range("I33").Value = 15
$ws.Range("J33").Value = 24
$ws.Range("K33").Value = -37.5
$ws.Range("L33").Value = -31.818181818181
